$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.361.38'

$ws.Range("D3").Value = '1.826.44'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.01'
$ws.Range("E5").Value = '  +0.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("E7").Value = '  -2.37%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3772'
$ws.Range("E8").Value = '  +0.91%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07459'

$ws.Range("E10").Value = '  +3.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.03'
$ws.Range("E11").Value = '  +0.12%  '

$ws.Range("D12").Value = '1.827.14'
$ws.Range("E12").Value = '  +0.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.755'
$ws.Range("E13").Value = '  +0.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.464'
$ws.Range("E14").Value = '  +2.07%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.07'
$ws.Range("E15").Value = '  +1.04%  '

$ws.Range("E16").Value = '  +0.54%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  -0.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008795'
$ws.Range("E18").Value = '  -0.62%  '

$ws.Range("E20").Value = '  +1.00%  '

$ws.Range("D21").Value = '27.370.36'
$ws.Range("E21").Value = '  +0.93%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.413'
$ws.Range("E22").Value = '  +4.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.99'
$ws.Range("E23").Value = '  -0.38%  '

$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.053.25'
$ws.Range("E24").Value = '  +0.36%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.963'
$ws.Range("E25").Value = '  -2.04%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.50'
$ws.Range("E26").Value = '  -0.19%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.327'
$ws.Range("E27").Value = '  +4.65%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.67'
$ws.Range("E28").Value = '  +0.91%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.416'
$ws.Range("E29").Value = '  +2.82%  '

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.01'
$ws.Range("E30").Value = '  +0.42%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08893'
$ws.Range("E31").Value = '  -0.05%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7923'
$ws.Range("E32").Value = '  +2.49%  '

$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.206'
$ws.Range("E33").Value = '  +0.92%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.602'
$ws.Range("E34").Value = '  +2.88%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.923'
$ws.Range("E35").Value = '  -1.71%  '

$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.0000'
$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.112'
$ws.Range("E37").Value = '  +0.71%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01989'
$ws.Range("E38").Value = '  +0.93%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05314'
$ws.Range("E39").Value = '  +0.49%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.315'
$ws.Range("E40").Value = '  +1.64%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5354'
$ws.Range("E41").Value = '  -0.85%  '

$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.879'
$ws.Range("E42").Value = '  -0.22%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1721'
$ws.Range("E43").Value = '  +0.35%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.300'
$ws.Range("E44").Value = '  +15.61%  '

$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.681'
$ws.Range("E45").Value = '  +0.57%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5111'
$ws.Range("E46").Value = '  -2.94%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.71'
$ws.Range("E47").Value = '  -0.21%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.699'
$ws.Range("E48").Value = '  +1.07%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.44'
$ws.Range("E49").Value = '  -0.64%  '

$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.000'
$ws.Range("E50").Value = '  -0.01%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06414'
$ws.Range("E51").Value = '  -1.21%  '
